# Atualização 06 e 07/07/2020
# Adds two new rows of COVID case data (06/07/2020 and 07/07/2020)
# to the end of the existing data table on Sheet1.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 100 : dia 99 / 06/07/2020 ---
$ws.Cells.Item(100, 1).Value  = 99
# Force the date column to be stored as plain text (matches the source
# data, which keeps dates as literal strings rather than date serials).
$ws.Cells.Item(100, 2).NumberFormat = "@"
$ws.Cells.Item(100, 2).Value  = "06/07/2020"
$ws.Cells.Item(100, 2).ClearFormats()
$ws.Cells.Item(100, 3).Value  = 919
$ws.Cells.Item(100, 4).Value  = 24
$ws.Cells.Item(100, 5).Value  = 124
$ws.Cells.Item(100, 6).Value  = "605,0510903"
$ws.Cells.Item(100, 7).Value  = "0,02611534276"
$ws.Cells.Item(100, 8).Value  = 767
$ws.Cells.Item(100, 9).Value  = 1197
$ws.Cells.Item(100, 10).Value = 2116
$ws.Cells.Item(100, 11).Value = 23
$ws.Cells.Item(100, 12).Value = 51
$ws.Cells.Item(100, 13).Value = 4
$ws.Cells.Item(100, 14).Value = 47
$ws.Cells.Item(100, 15).Value = 66
$ws.Cells.Item(100, 16).Value = 21
$ws.Cells.Item(100, 17).Value = 18
$ws.Cells.Item(100, 18).Value = 15

# --- Row 101 : dia 100 / 07/07/2020 ---
$ws.Cells.Item(101, 1).Value  = 100
$ws.Cells.Item(101, 2).NumberFormat = "@"
$ws.Cells.Item(101, 2).Value  = "07/07/2020"
$ws.Cells.Item(101, 2).ClearFormats()
$ws.Cells.Item(101, 3).Value  = 947
$ws.Cells.Item(101, 4).Value  = 24
$ws.Cells.Item(101, 5).Value  = 119
$ws.Cells.Item(101, 6).Value  = "623,4857263"
$ws.Cells.Item(101, 7).Value  = "0,02534318902"
$ws.Cells.Item(101, 8).Value  = 800
$ws.Cells.Item(101, 9).Value  = 1219
$ws.Cells.Item(101, 10).Value = 2166
$ws.Cells.Item(101, 11).Value = 28
$ws.Cells.Item(101, 12).Value = 54
$ws.Cells.Item(101, 13).Value = 4
$ws.Cells.Item(101, 14).Value = 50
$ws.Cells.Item(101, 15).Value = 50
$ws.Cells.Item(101, 16).Value = 26
$ws.Cells.Item(101, 17).Value = 17
$ws.Cells.Item(101, 18).Value = 15
